$wb = $excel.ActiveWorkbook

# Rename the second sheet from "sprint0" to "sprint3"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "sprint3"

# Make the renamed sheet the active/selected tab (was Backlog before)
$ws2.Select()
